$wb = $excel.ActiveWorkbook
$nl = [char]10

# ---------------------------------------------------------------------------
# Sheet 2: "Dizionario_Relazioni" -- merge the old "Prenota" + "Acquista"
# relationships into a single "Prenota/Acquista" relationship row, and
# delete the now-redundant "Acquista" row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B6").Value = "Prenota/" + $nl + "Acquista"
$ws2.Range("B6").WrapText = $true

$ws2.Range("D6").Value = "Itinerario" + $nl + "Passeggero"
$ws2.Range("D6").WrapText = $true

$ws2.Range("C6").Value = "Un passeggero può prenotare/acquistare più itinerari"
$ws2.Range("C6").WrapText = $false

$ws2.Rows("6").RowHeight = 31.2

# Remove the obsolete "Acquista" relationship row (old row 8).
$ws2.Rows("8").Delete()

# ---------------------------------------------------------------------------
# Sheet 1: "Dizionario_Entità" -- update the Attributi text for Dati_Pagamento
# (AccountID -> IdAccount) and move the selection.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D3").Value = "IdDatiPagamento, Tipo, Nome, Cognome, Paese, Via, Città, CAP, NumeroCarta, DataValidita, IdAccount"

$ws2.Activate()
$ws2.Range("D17").Select()

$ws1.Activate()
$ws1.Range("D3").Select()
